# Add a new "Logged Work" test case as a second worksheet (Sheet2),
# following the existing "Sheet1" test-case list layout/format.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Sheet1's selection (no data changes on Sheet1 itself) ---
$null = $ws1.Range("A1:J4").Select()

# --- Create Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Header row (reuses the same column headers as Sheet1)
$ws2.Range("A1").Value = "Test Case ID"
$ws2.Range("B1").Value = "Title / Summary"
$ws2.Range("C1").Value = "Description"
$ws2.Range("D1").Value = "Precondition"
$ws2.Range("E1").Value = "Test Steps"
$ws2.Range("F1").Value = "Expected Result"
$ws2.Range("G1").Value = "Test Data"
$ws2.Range("H1").Value = "Priority"
$ws2.Range("I1").Value = "Status"
$ws2.Range("J1").Value = "Module"
$ws2.Range("A1:J1").Font.Bold = $true

# New test case: TC_LogWorkIn_001 (write order matches shared-string creation order)
$ws2.Range("A2").Value = "TC_LogWorkIn_001"
$ws2.Range("C2").Value = "User Loged Work (in)"
$ws2.Range("B2").Value = "Logged Work"
$ws2.Range("D2").Value = "User on dashboard page"
$ws2.Range("F2").Value = "Success alert pop up"
$ws2.Range("I2").Value = "On Progress"

$ws2.Range("E2").Value = '1. Click Menu "Time"'
$ws2.Range("E3").Value = '2. Click Top Bar "Attendance"'
$ws2.Range("E4").Value = "3. Click Punch In/Out"
$ws2.Range("E5").Value = "4. Select Desired Date"
$ws2.Range("E6").Value = "5. Select Desired Time"
$ws2.Range("E7").Value = "6. Click In"

$ws2.Range("H2").Value = "High"

# Merge the per-test-case columns down across the 6 step rows (2-7)
$ws2.Range("A2:A7").Merge()
$ws2.Range("B2:B7").Merge()
$ws2.Range("C2:C7").Merge()
$ws2.Range("D2:D7").Merge()
$ws2.Range("F2:F7").Merge()
$ws2.Range("H2:H7").Merge()
$ws2.Range("I2:I7").Merge()
$ws2.Range("J2:J7").Merge()

# Title/Summary column is centered for this test case
$ws2.Range("B2:B7").HorizontalAlignment = -4108

# Touch the remaining merged ranges with a no-op border so every row of the
# merge keeps a (formatted) cell record, matching the rest of the workbook's
# merged blocks (e.g. Sheet1 A8:A10).
$ws2.Range("A2:A7").Borders.LineStyle = -4142
$ws2.Range("C2:C7").Borders.LineStyle = -4142
$ws2.Range("D2:D7").Borders.LineStyle = -4142
$ws2.Range("F2:F7").Borders.LineStyle = -4142
$ws2.Range("H2:H7").Borders.LineStyle = -4142
$ws2.Range("I2:I7").Borders.LineStyle = -4142
$ws2.Range("J2:J7").Borders.LineStyle = -4142

# Column widths (calibrated to match the target OOXML "width" values)
$ws2.Columns.Item(1).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 18.5
$ws2.Columns.Item(4).ColumnWidth = 23.666666666666668
$ws2.Columns.Item(5).ColumnWidth = 28.666666666666668
$ws2.Columns.Item(6).ColumnWidth = 32.166666666666664
$ws2.Columns.Item(7).ColumnWidth = 24.833333333333332
$ws2.Columns.Item(10).ColumnWidth = 12.333333333333334

# Activate Sheet2 and set its selection
$ws2.Activate()
$null = $ws2.Range("F21").Select()
